# "Moved ramps to unit__nodes"
#
# The ramp_cost / ramp_method / ramp_speed_down / ramp_speed_up columns
# (T:W) are removed from the unit_c sheet and the same four columns
# (with their header text and formatting) are added to the
# unit_sourceNode_c and unit_sinkNode_c sheets at F:I.

$wb = $excel.ActiveWorkbook

$unitC        = $wb.Worksheets.Item("unit_c")
$sourceNodeC  = $wb.Worksheets.Item("unit_sourceNode_c")
$sinkNodeC    = $wb.Worksheets.Item("unit_sinkNode_c")

# Copy the ramp columns (header description row 1, header name row 2)
# from unit_c, including their formatting, to F1 on the two node sheets.
$unitC.Range("T1:W2").Copy()
$sourceNodeC.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$sourceNodeC.Range("F1").PasteSpecial(-4163)  # xlPasteValues

$unitC.Range("T1:W2").Copy()
$sinkNodeC.Range("F1").PasteSpecial(-4122)    # xlPasteFormats
$sinkNodeC.Range("F1").PasteSpecial(-4163)    # xlPasteValues

$excel.CutCopyMode = 0

# Now that the ramp data lives on the node sheets, remove the original
# columns from unit_c. The remaining columns (retire_max_total,
# startup_cost, startup_method, variable_cost, virtual_unitsize) shift
# left to take their place.
$unitC.Range("T1:W2").EntireColumn.Delete()

# Reproduce the selections left behind on the two node sheets.
$sourceNodeC.Activate()
$sourceNodeC.Range("F1:I2").Select()

$sinkNodeC.Activate()
$sinkNodeC.Range("G1").Select()

# unit_c becomes the active/selected sheet (previously node_c was).
$unitC.Activate()
$unitC.Range("A1").Select()
